$d = $word.ActiveDocument

# 1. Replace the title text: remove ", well," from the title run
$d.Content.Find.Execute("The OECD Fragility Clusters are, well, Fragile", $true, $false, $false, $false, $false, $true, 1, $false, "The OECD Fragility Clusters are Fragile", 2)

Write-Host "Done"
